$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "year" column (D) values for the relevant rows.
# Rows 2-6   -> 2023
# Rows 7-16  -> 2024
# Rows 17-28 -> 2025

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = 2023
}

for ($r = 7; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = 2024
}

for ($r = 17; $r -le 28; $r++) {
    $ws.Cells.Item($r, 4).Value = 2025
}
